# "Disconnected coke from pellet production":
#   The 'connections' sheet had a row (row 5) wiring the fossil-fuel/coke
#   outflow into the pellets chain via simple_coke/simple_pellets. That
#   connection row is removed entirely, and every row below it shifts up
#   by one (old row 6 -> new row 5, ... old row 18 -> new row 17), leaving
#   the former last row (old row 19, blank placeholder with just styles)
#   as the new last row (18).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Row 5 was: pellets | simple_pellets | inflow | fossil fuel | coke |
#            outflows | simple_coke | coke
# Deleting the whole row shifts everything below up by one, which is
# exactly the change shown in the diff (rows 6-19 -> 5-18).
$ws.Rows(5).Delete() | Out-Null

# The edited workbook was left with the "connections" sheet active/on
# screen (instead of "CO2"), with the cursor resting on B22.
$ws.Activate()
$ws.Range("B22").Select() | Out-Null
